$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Note" column (Q) for rows 2-38 to reflect that the bank
# location is now operating normally instead of awaiting a bank update.
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 17).Value = "เปิดทำการปกติ"
}

# Move the active selection from the logo column (P) to the note column (Q),
# matching the updated selection saved in the workbook.
$ws.Range("Q2:Q38").Select()
